$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing value in A16 (tiny float correction)
$ws.Range("A16").Value = 45864.83358131944

# Add the new row 17 data
$ws.Range("A17").Value = 45864.87532446441
$ws.Range("B17").Value = 2025
$ws.Range("C17").Value = 30
$ws.Range("D17").Value = 13.61
$ws.Range("E17").Value = 89.84999999999999
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 4.19
$ws.Range("H17").Value = "E"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "21:00:28"

# Copy the style from A16 (date style) to A17
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122) # xlPasteFormats
